# Update Lrp1 (App-Lrp1) LR-pair sheet with refreshed TPM-derived values.
#
# The underlying per-cluster averages (ligand "G"/"H" and receptor "M"/"N")
# change for the "ECs" cluster because of the new TPM input; every other
# column on the sheet (detection-rate specificities and edge weights) is a
# deterministic function of those per-cluster values, so we recompute them
# here instead of hard-coding every cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 10

# --- 1. New ligand (sending-cluster) and receptor (target-cluster) base values ---
# Only the "ECs" cluster's ligand/receptor expression values changed.
$newG = @{ "ECs" = 97.57717366666668 }
$newH = @{ "ECs" = 292.731521 }
$newM = @{ "ECs" = 3.456265333333333 }
$newN = @{ "ECs" = 10.368796 }

# --- 2. Read current per-row cluster labels + G/H/M/N, applying overrides ---
$sendOf = @{}
$targOf = @{}
$Gof = @{}
$Hof = @{}
$Mof = @{}
$Nof = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()   # A: Sending cluster
    $d = $ws.Cells.Item($r, 4).Value()   # D: Target cluster
    $sendOf[$r] = $a
    $targOf[$r] = $d

    if ($newG.ContainsKey($a)) {
        $g = $newG[$a]
    } else {
        $g = $ws.Cells.Item($r, 7).Value()
    }
    if ($newH.ContainsKey($a)) {
        $h = $newH[$a]
    } else {
        $h = $ws.Cells.Item($r, 8).Value()
    }
    if ($newM.ContainsKey($d)) {
        $m = $newM[$d]
    } else {
        $m = $ws.Cells.Item($r, 13).Value()
    }
    if ($newN.ContainsKey($d)) {
        $n = $newN[$d]
    } else {
        $n = $ws.Cells.Item($r, 14).Value()
    }

    $Gof[$a] = $g
    $Hof[$a] = $h
    $Mof[$d] = $m
    $Nof[$d] = $n
}

# --- 3. Totals across the distinct clusters (for specificity ratios) ---
$sumG = 0
$sumH = 0
foreach ($k in $Gof.Keys) { $sumG += $Gof[$k] }
foreach ($k in $Hof.Keys) { $sumH += $Hof[$k] }

$sumM = 0
$sumN = 0
foreach ($k in $Mof.Keys) { $sumM += $Mof[$k] }
foreach ($k in $Nof.Keys) { $sumN += $Nof[$k] }

# --- 4. Write G/H/M/N, then the edge Q/R values per row ---
$Qof = @{}
$Rof = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $a = $sendOf[$r]
    $d = $targOf[$r]

    $g = $Gof[$a]
    $h = $Hof[$a]
    $m = $Mof[$d]
    $n = $Nof[$d]

    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n

    $i = $g / $sumG
    $j = $h / $sumH
    $o = $m / $sumM
    $p = $n / $sumN

    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p

    $q = $g * $m
    $rr = $h * $n
    $Qof[$r] = $q
    $Rof[$r] = $rr

    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rr
}

# --- 5. Edge-weight specificities (depend on the total across all rows) ---
$sumQ = 0
$sumR = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sumQ += $Qof[$r]
    $sumR += $Rof[$r]
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $s = $Qof[$r] / $sumQ
    $t = $Rof[$r] / $sumR
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}
